$wb = $excel.ActiveWorkbook

$oldGuid = "549ffc53-08a6-4c9c-ae30-decbfd0778ba"
$newGuid = "b2ab9d98-49ce-4762-99de-f481f6e96c7e"

$oldHash = "5570f495a7fb63e3e7e127377d7369e54ce02485"
$newHash = "6445cc3499957e601b12462798f436716c898ee3"

$newMdName = "$newGuid.md"
$newZhName = "$newGuid.$newHash.zh-cn.xlf"
$newDeName = "$newGuid.$newHash.de-de.xlf"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/2461f0286ed73a70b2363edbdaf2e3cddb9e2b0b/e2e/$oldGuid.md"
$zhTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf58cdf8d4fcdd694099e5e2a6c99248837f0ebf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b9580d709573fbe2ab614080d7b657820617500/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-56-21 04:56:02"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMdName) | Out-Null

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhName
$wsZh.Range("E2").Value = "2016-03-21 04:55:59"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdTarget, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdTarget, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhTarget, "", "", $newZhName) | Out-Null

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeName
$wsDe.Range("E2").Value = "2016-03-21 04:56:02"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdTarget, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdTarget, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deTarget, "", "", $newDeName) | Out-Null
